$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at P:Q, shifting old P,Q,R,S -> R,S,T,U
$ws.Range("P1:Q1").EntireColumn.Insert()

# New header labels
$ws.Range("P1").Value = "ligand_conc"
$ws.Range("Q1").Value = "time_step"
$ws.Range("R1").Value = "comments"
$ws.Range("S1").Value = "assay_date"

# Column widths for the newly inserted columns
$ws.Range("P1").EntireColumn.ColumnWidth = 17
$ws.Range("Q1").EntireColumn.ColumnWidth = 17

# Column M width changes (was bestFit shared with N, now its own custom width)
$ws.Range("M1").EntireColumn.ColumnWidth = 22.54296875

# New data cells on rows 4 and 5 to match style of column O
$ws.Range("P4").Value = $null
$ws.Range("Q4").Value = $null
$ws.Range("P5").Value = $null
$ws.Range("Q5").Value = $null

# Sheet view changes
$ws.Application.ActiveWindow.Zoom = 83
$ws.Range("Q7").Select()
